# The commit swaps the two theme parts of the deck: the slide master's
# theme (ppt/theme/theme1.xml, "Integral") ends up holding the color
# scheme that used to live in ppt/theme/theme2.xml ("Office Theme"), and
# vice-versa. The exposed PowerPoint object model only surfaces a single
# editable Theme (the one bound to the presentation's slide master /
# Designs(1)), reachable through ThemeColorScheme, so we recolor that
# theme's 12 color slots to the target "Office Theme" palette.
#
# COM RGB() values are 0x00BBGGRR (R + G*256 + B*65536), matching the
# hex srgbClr values that end up in the OOXML in R,G,B order.

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette: the "Office Theme" clrScheme that currently lives in
# ppt/theme/theme2.xml, in dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# order (ThemeColorScheme.Item(1..12)).
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToRgb($officeThemeColors[$i - 1])
}
